# Referee.xlsx - "Tested with referee config."
# Replace the "1000 100 50" / "1000 200 50" parameter strings used by the
# referee LED/buzzer FLASH actions with a single "1000 200 200" value, fill
# in the previously empty F16 parameter cell, and change the B23 / D19
# counters to be stored as text ("4" / "3") instead of numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- F10, F11 and F16 need their cell format copied from F2 (no fill,
#     right aligned text) before the value is written, since in the
#     original workbook they used the shaded-row format (style shared with
#     G10 / G11, which must stay untouched).
$ws.Range("F2").Copy()
$ws.Range("F10").PasteSpecial(-4122)
$ws.Range("F11").PasteSpecial(-4122)
$ws.Range("F16").PasteSpecial(-4122)

# --- B23 and D19 change from a plain number to a text value.
$ws.Range("B23").Value = "4"
$ws.Range("D19").Value = "3"

# --- Update the "parameters" column (F) for every FLASH row to the new
#     single value "1000 200 200".
$ws.Range("F2").Value = "1000 200 200"
$ws.Range("F3").Value = "1000 200 200"
$ws.Range("F7").Value = "1000 200 200"
$ws.Range("F10").Value = "1000 200 200"
$ws.Range("F11").Value = "1000 200 200"
$ws.Range("F16").Value = "1000 200 200"
$ws.Range("F18").Value = "1000 200 200"
$ws.Range("F19").Value = "1000 200 200"
$ws.Range("F23").Value = "1000 200 200"

# --- Leave the selection where the editing ended up.
$ws.Range("F23").Select()
